$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The BOM no longer lists a standalone "N-DAP" placeholder part (row 2:
# Comment="N-DAP", Description="", Footprint="N-DAP", Designator="B1").
# Removing it shifts every subsequent part row up by one.
$ws.Range("A2:E2").EntireRow.Delete()

# Column widths were re-tuned slightly (Eurocircuits assembly output run).
$ws.Columns.Item(1).ColumnWidth = 11
$ws.Columns.Item(2).ColumnWidth = 23.166666666666668
$ws.Columns.Item(3).ColumnWidth = 40.833333333333336
$ws.Columns.Item(4).ColumnWidth = 22.166666666666668
$ws.Columns.Item(5).ColumnWidth = 69.83333333333333
